# DOMA-1872: remove multi-tariff values (2nd/3rd readings) for non-electricity
# meters in meter-import-example.xlsx. Only meters of type "ЭЛ" (electricity)
# keep a tariff count of 3 and the extra reading values; every other meter
# type (ГВС, ХВС, ТЕПЛО, ГАЗ) is reduced to a single tariff with readings
# 2/3 cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the used range / last row of data.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 1 }

for ($r = 2; $r -le $lastRow; $r++) {
    $meterType = $ws.Cells.Item($r, 4).Value()  # column D = "Тип счетчика"
    if ($meterType -ne "ЭЛ") {
        $ws.Cells.Item($r, 6).Value = 1        # column F = "Кол-во тарифов"
        $ws.Cells.Item($r, 8).ClearContents()  # column H = "Показание 2"
        $ws.Cells.Item($r, 9).ClearContents()  # column I = "Показание 3"
    }
}

# Row 10's "№ счетчика" (column E) was re-numbered from 22 to 33 to keep
# meter numbers unique within the example data.
$ws.Cells.Item(10, 5).Value = 33
